# Publish refresh: bump as_of_utc timestamps and update the referee stat
# rows that changed between snapshots (2025-11-25 07:09:29 -> 2025-11-26 03:06:02).

$wb = $excel.ActiveWorkbook

$oldStamp = "2025-11-25 07:09:29"
$newStamp = "2025-11-26 03:06:02"

# ---------------------------------------------------------------------------
# Sheet "Главные" (main officials)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Главные")

# Refresh the as_of_utc column (AA) for every data row (2..26).
for ($r = 2; $r -le 26; $r++) {
    $ws2.Cells.Item($r, 27).Value = $newStamp
}

# Row 8 - Gamaley Evgeniy
$ws2.Range("C8").Value = 26
$ws2.Range("D8").Value = 458
$ws2.Range("E8").Value = 231
$ws2.Range("F8").Value = 227
$ws2.Range("G8").Value = 17.62
$ws2.Range("H8").Value = 8.880000000000001
$ws2.Range("I8").Value = 8.73
$ws2.Range("J8").Value = 108
$ws2.Range("K8").Value = 106

# Row 9 - Gashilov Viktor
$ws2.Range("C9").Value = 29
$ws2.Range("D9").Value = 446
$ws2.Range("E9").Value = 233
$ws2.Range("F9").Value = 213
$ws2.Range("G9").Value = 15.38
$ws2.Range("H9").Value = 8.029999999999999
$ws2.Range("I9").Value = 7.34
$ws2.Range("J9").Value = 114
$ws2.Range("K9").Value = 104
$ws2.Range("V9").Value = 16
$ws2.Range("W9").Value = 30

# Row 15 - Lazarev Gleb
$ws2.Range("C15").Value = 19
$ws2.Range("D15").Value = 348
$ws2.Range("E15").Value = 165
$ws2.Range("F15").Value = 183
$ws2.Range("G15").Value = 18.32
$ws2.Range("H15").Value = 8.68
$ws2.Range("I15").Value = 9.630000000000001
$ws2.Range("J15").Value = 60
$ws2.Range("K15").Value = 79
$ws2.Range("L15").Value = 5
$ws2.Range("M15").Value = 5

# Row 16 - Morozov Sergey
$ws2.Range("C16").Value = 28
$ws2.Range("D16").Value = 511
$ws2.Range("E16").Value = 254
$ws2.Range("F16").Value = 257
$ws2.Range("G16").Value = 18.25
$ws2.Range("H16").Value = 9.07
$ws2.Range("I16").Value = 9.18
$ws2.Range("J16").Value = 97
$ws2.Range("K16").Value = 96

# Row 20 - Oskirko Yuriy
$ws2.Range("C20").Value = 27
$ws2.Range("D20").Value = 475
$ws2.Range("E20").Value = 205
$ws2.Range("F20").Value = 270
$ws2.Range("G20").Value = 17.59
$ws2.Range("H20").Value = 7.59
$ws2.Range("I20").Value = 10
$ws2.Range("J20").Value = 95
$ws2.Range("K20").Value = 100
$ws2.Range("L20").Value = 3
$ws2.Range("M20").Value = 8

# ---------------------------------------------------------------------------
# Sheet "Линейные" (linesmen)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Линейные")

# Refresh the as_of_utc column (AA) for every data row (2..26).
for ($r = 2; $r -le 26; $r++) {
    $ws3.Cells.Item($r, 27).Value = $newStamp
}

# Row 3 - Bersenyov Maksim
$ws3.Range("C3").Value = 28
$ws3.Range("D3").Value = 405
$ws3.Range("E3").Value = 209
$ws3.Range("F3").Value = 196
$ws3.Range("G3").Value = 14.46
$ws3.Range("H3").Value = 7.46
$ws3.Range("I3").Value = 7
$ws3.Range("J3").Value = 102
$ws3.Range("K3").Value = 83

# Row 9 - Golovlyov Dmitriy
$ws3.Range("C9").Value = 27
$ws3.Range("D9").Value = 511
$ws3.Range("E9").Value = 220
$ws3.Range("F9").Value = 291
$ws3.Range("G9").Value = 18.93
$ws3.Range("H9").Value = 8.15
$ws3.Range("I9").Value = 10.78
$ws3.Range("J9").Value = 95
$ws3.Range("K9").Value = 118
$ws3.Range("L9").Value = 4
$ws3.Range("M9").Value = 5

# Row 18 - Novikov Nikita
$ws3.Range("C18").Value = 29
$ws3.Range("D18").Value = 475
$ws3.Range("E18").Value = 226
$ws3.Range("F18").Value = 249
$ws3.Range("G18").Value = 16.38
$ws3.Range("H18").Value = 7.79
$ws3.Range("I18").Value = 8.59
$ws3.Range("J18").Value = 108
$ws3.Range("K18").Value = 107
$ws3.Range("V18").Value = 6
$ws3.Range("W18").Value = 22

# Row 19 - Parikov Yaroslav
$ws3.Range("C19").Value = 26
$ws3.Range("D19").Value = 435
$ws3.Range("E19").Value = 208
$ws3.Range("F19").Value = 227
$ws3.Range("G19").Value = 16.73
$ws3.Range("H19").Value = 8
$ws3.Range("I19").Value = 8.73
$ws3.Range("J19").Value = 99
$ws3.Range("K19").Value = 101
$ws3.Range("V19").Value = 10
$ws3.Range("W19").Value = 14

# Row 26 - Slavikovskiy Roman
$ws3.Range("C26").Value = 26
$ws3.Range("D26").Value = 540
$ws3.Range("E26").Value = 230
$ws3.Range("F26").Value = 310
$ws3.Range("G26").Value = 20.77
$ws3.Range("H26").Value = 8.85
$ws3.Range("I26").Value = 11.92
$ws3.Range("J26").Value = 90
$ws3.Range("K26").Value = 95
$ws3.Range("L26").Value = 8
$ws3.Range("M26").Value = 10
